$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# zh-cn sheet (row 7 = aa830a94-b24b-4fc6-b011-78bc40ae6624.md handback report)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7b6b4d8cb8c5122b80f791bb0bea913fda4d1fce/e2e/aa830a94-b24b-4fc6-b011-78bc40ae6624.md",
    "",
    "",
    "aa830a94-b24b-4fc6-b011-78bc40ae6624.md"
) | Out-Null

$wsZh.Range("J7").Value = "aa830a94-b24b-4fc6-b011-78bc40ae6624.a8b9c96d33a278e161c50412fa05e022e2df6f51.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-19 02:52:30"
$wsZh.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6edf89d2815e63da257c2d9c73d7dbdaf675d830/e2e/aa830a94-b24b-4fc6-b011-78bc40ae6624.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6e7b91e39bc51e81a2617480081cc1dedb84aef3/e2e/aa830a94-b24b-4fc6-b011-78bc40ae6624.md."

# ---------------------------------------------------------------------------
# de-de sheet (row 7 = aa830a94-b24b-4fc6-b011-78bc40ae6624.md handback report)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5244bc44b26783f7893a2d147469cced00773458/e2e/aa830a94-b24b-4fc6-b011-78bc40ae6624.md",
    "",
    "",
    "aa830a94-b24b-4fc6-b011-78bc40ae6624.md"
) | Out-Null

$wsDe.Range("J7").Value = "aa830a94-b24b-4fc6-b011-78bc40ae6624.a8b9c96d33a278e161c50412fa05e022e2df6f51.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-19 02:52:37"
$wsDe.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6edf89d2815e63da257c2d9c73d7dbdaf675d830/e2e/aa830a94-b24b-4fc6-b011-78bc40ae6624.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6e7b91e39bc51e81a2617480081cc1dedb84aef3/e2e/aa830a94-b24b-4fc6-b011-78bc40ae6624.md."
